$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "66.158.94"
$ws.Range("E2").Value = "  +1.22%  "
# Row 3
$ws.Range("D3").Value = "3.562.77"
$ws.Range("E3").Value = "  +5.19%  "
# Row 4
$ws.Range("E4").Value = "  +0.13%  "
# Row 5
$ws.Range("D5").Value = "607.34"
$ws.Range("E5").Value = "  +2.25%  "
# Row 6
$ws.Range("D6").Value = "144.84"
$ws.Range("E6").Value = "  +2.59%  "
# Row 7
$ws.Range("D7").Value = "3.562.01"
$ws.Range("E7").Value = "  +5.16%  "
# Row 8
$ws.Range("E8").Value = "  +0.10%  "
# Row 9
$ws.Range("E9").Value = "  +4.29%  "
# Row 10
$ws.Range("E10").Value = "  +2.58%  "
# Row 11
$ws.Range("D11").Value = "7.95"
$ws.Range("E11").Value = "  +0.65%  "
# Row 12
$ws.Range("E12").Value = "  +1.57%  "
# Row 13
$ws.Range("D13").Value = "4.168.30"
$ws.Range("E13").Value = "  +5.30%  "
# Row 14
$ws.Range("E14").Value = "  +4.43%  "
# Row 15
$ws.Range("D15").Value = "30.18"
$ws.Range("E15").Value = "  +2.00%  "
# Row 16
$ws.Range("D16").Value = "3.563.81"
$ws.Range("E16").Value = "  +5.40%  "
# Row 17
$ws.Range("D17").Value = "66.277.43"
$ws.Range("E17").Value = "  +1.38%  "
# Row 18
$ws.Range("E18").Value = "  -0.55%  "
# Row 19
$ws.Range("D19").Value = "11.47"
$ws.Range("E19").Value = "  +10.38%  "
# Row 20
$ws.Range("E20").Value = "  +1.87%  "
# Row 21
$ws.Range("D21").Value = "14.92"
$ws.Range("E21").Value = "  +1.91%  "
# Row 22
$ws.Range("D22").Value = "431.28"
$ws.Range("E22").Value = "  +4.15%  "
# Row 23
$ws.Range("E23").Value = "  +5.71%  "
# Row 24
$ws.Range("D24").Value = "78.83"
$ws.Range("E24").Value = "  +2.57%  "
# Row 25
$ws.Range("D25").Value = "3.705.07"
$ws.Range("E25").Value = "  +5.33%  "
# Row 26
$ws.Range("E26").Value = "  -0.07%  "
# Row 27
$ws.Range("E27").Value = "  +9.72%  "
# Row 28
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +4.89%  "
# Row 29
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  +3.64%  "
# Row 30
$ws.Range("D30").Value = "9.12"
$ws.Range("E30").Value = "  -0.90%  "
# Row 31
$ws.Range("E31").Value = "  -0.25%  "
# Row 32
$ws.Range("E32").Value = "  +2.13%  "
# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "25.51"
$ws.Range("E33").Value = "  +5.08%  "
# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.158"
$ws.Range("E34").Value = "  -1.16%  "
# Row 35
$ws.Range("D35").Value = "3.558.64"
$ws.Range("E35").Value = "  +5.21%  "
# Row 36
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "
# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  +4.90%  "
# Row 38
$ws.Range("D38").Value = "7.92"
$ws.Range("E38").Value = "  +5.60%  "
# Row 39
$ws.Range("E39").Value = "  +2.32%  "
# Row 40
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.16%  "
# Row 41
$ws.Range("D41").Value = "170.09"
$ws.Range("E41").Value = "  +2.44%  "
# Row 42
$ws.Range("D42").Value = "0.0855"
$ws.Range("E42").Value = "  +0.24%  "
# Row 43
$ws.Range("E43").Value = "  +4.08%  "
# Row 44
$ws.Range("E44").Value = "  +3.38%  "
# Row 45
$ws.Range("E45").Value = "  +1.65%  "
# Row 46
$ws.Range("D46").Value = "46.16"
$ws.Range("E46").Value = "  +1.99%  "
# Row 47
$ws.Range("E47").Value = "  +3.54%  "
# Row 48
$ws.Range("D48").Value = "25.92"
$ws.Range("E48").Value = "  -2.36%  "
# Row 49
$ws.Range("E49").Value = "  +5.73%  "
# Row 50
$ws.Range("E50").Value = "  +1.83%  "
# Row 51
$ws.Range("D51").Value = "23.50"
$ws.Range("E51").Value = "  +16.01%  "
